$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the dialogue table values ---
$ws.Range('A1').Value = 'DIALOGUE'
$ws.Range('B1').Value = 'english'
$ws.Range('C1').Value = 'czech'
$ws.Range('D1').Value = 'korean'
$ws.Range('E1').Value = 'japanese'

$ws.Range('A2').Value = 'headingFont'
$ws.Range('B2').Value = 'heading_EN.tff'
$ws.Range('C2').Value = 'heading_EN.tff'
$ws.Range('D2').Value = 'heading_KOR.tff'
$ws.Range('E2').Value = 'heading_JAP.otf'

$ws.Range('A3').Value = 'textFont'
$ws.Range('B3').Value = 'text_EN.tff'
$ws.Range('C3').Value = 'text_CZ.tff'
$ws.Range('D3').Value = 'text_EN.tff'
$ws.Range('E3').Value = 'text_JAP.otf'

$ws.Range('A4').Value = 'languageDisplayName'
$ws.Range('B4').Value = 'English'
$ws.Range('C4').Value = 'Čeština'
$ws.Range('D4').Value = '한국인'
$ws.Range('E4').Value = '日本語'
$ws.Rows.Item(4).RowHeight = 42

$ws.Range('A5').Value = 'appName'
$ws.Range('B5').Value = 'Shooting Stars'
$ws.Range('C5').Value = 'Padající Hvězdy'
$ws.Range('D5').Value = '별을 쏘다'
$ws.Range('E5').Value = '流れ星撃ち'
$ws.Rows.Item(5).RowHeight = 42

$ws.Range('A6').Value = 'menuSubText'
$ws.Range('B6').Value = 'Press any key to continue'
$ws.Range('C6').Value = 'Stiskněte libovolnou klávesu pro pokračování'
$ws.Range('D6').Value = '아무 키나 누르세요 계속하려면'
$ws.Range('E6').Value = '続行するには任意のキーを押してください'
$ws.Rows.Item(6).RowHeight = 42

$ws.Range('A7').Value = 'score'
$ws.Range('B7').Value = 'Score'
$ws.Range('C7').Value = 'Skóre'
$ws.Range('D7').Value = '점수'
$ws.Range('E7').Value = 'スコア'
$ws.Rows.Item(7).RowHeight = 42

$ws.Range('A8').Value = 'topscore'
$ws.Range('B8').Value = 'Top score'
$ws.Range('C8').Value = 'Nejlepší skóre'
$ws.Range('D8').Value = '최고 점수'
$ws.Range('E8').Value = 'トプスコア'
$ws.Rows.Item(8).RowHeight = 42

$ws.Range('A9').Value = 'timeLeft'
$ws.Range('B9').Value = 'Time left'
$ws.Range('C9').Value = 'Zbývající čas'
$ws.Range('D9').Value = '남은 시간'
$ws.Range('E9').Value = '残り時間'
$ws.Rows.Item(9).RowHeight = 42

$ws.Range('A10').Value = 'gameOver'
$ws.Range('B10').Value = 'GAME OVER'
$ws.Range('C10').Value = 'KONEC HRY'
$ws.Range('D10').Value = '게임 오버'
$ws.Range('E10').Value = 'ゲームオーバー'
$ws.Rows.Item(10).RowHeight = 42

$ws.Range('A11').Value = 'gameOverSubtext'
$ws.Range('B11').Value = 'Press "R" to restart'
$ws.Range('C11').Value = 'Stiskněte "R" pro restart'
$ws.Range('D11').Value = '다시 시작하려면 "R"을 누르세요'
$ws.Range('E11').Value = '再起動するには "R" を押してください'
$ws.Rows.Item(11).RowHeight = 42

$ws.Range('A12').Value = 'pause'
$ws.Range('B12').Value = 'Pause'
$ws.Range('C12').Value = 'Pozastaveno'
$ws.Range('D12').Value = '일시 정지'
$ws.Range('E12').Value = '一時停止'
$ws.Rows.Item(12).RowHeight = 42

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 50.65
$ws.Columns.Item(3).ColumnWidth = 64.15
$ws.Columns.Item(4).ColumnWidth = 59.35
$ws.Columns.Item(5).ColumnWidth = 73.0

# --- Freeze panes on column A, select B3 ---
$ws.Range('B1').Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range('B3').Select() | Out-Null
